$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AZ3").Value = "TN"
$ws.Range("AV5").Value = "TN"
$ws.Range("AR6").Value = "TN"
$ws.Range("AV7").Value = "TN"
$ws.Range("AZ9").Value = "FP"
$ws.Range("AZ10").Value = "FP"
$ws.Range("AZ11").Value = "TN"
$ws.Range("AR12").Value = "TN"
$ws.Range("AV12").Value = "FP"
$ws.Range("AQ13").Value = "FP"
$ws.Range("AR15").Value = "TN"
$ws.Range("AV17").Value = "TN"
$ws.Range("AR19").Value = "TN"
$ws.Range("AZ20").Value = "FP"
$ws.Range("AR22").Value = "TN"
$ws.Range("AQ23").Value = "FP"
$ws.Range("AR24").Value = "TN"
$ws.Range("AZ24").Value = "TN"
$ws.Range("AQ31").Value = "TN"
$ws.Range("AV31").Value = "TN"
$ws.Range("AQ32").Value = "FP"
$ws.Range("AR32").Value = "TN"
$ws.Range("AQ33").Value = "TN"
$ws.Range("AQ37").Value = "TN"
$ws.Range("AZ37").Value = "FP"
$ws.Range("AR40").Value = "FP"
$ws.Range("AQ41").Value = "FP"
$ws.Range("AR42").Value = "FP"
$ws.Range("AR43").Value = "TN"
$ws.Range("AV43").Value = "TN"
$ws.Range("AR44").Value = "TN"
$ws.Range("AZ48").Value = "FP"
$ws.Range("AZ52").Value = "FP"
$ws.Range("AZ53").Value = "FP"
$ws.Range("AZ54").Value = "FP"
$ws.Range("AZ56").Value = "TN"
$ws.Range("AZ58").Value = "FP"
$ws.Range("AV62").Value = "FN"
$ws.Range("AR63").Value = "FP"
$ws.Range("AZ64").Value = "TN"
$ws.Range("AV71").Value = "TN"
$ws.Range("AR72").Value = "FN"
$ws.Range("AV73").Value = "TN"
$ws.Range("AZ74").Value = "FN"
$ws.Range("AR75").Value = "FP"
$ws.Range("AZ75").Value = "TN"
$ws.Range("AZ76").Value = "TN"
$ws.Range("AR78").Value = "TN"
$ws.Range("AV78").Value = "TN"
$ws.Range("AQ79").Value = "FP"
$ws.Range("AZ80").Value = "FP"
$ws.Range("AQ86").Value = "TN"
$ws.Range("AV87").Value = "TN"
$ws.Range("AQ88").Value = "FP"
$ws.Range("AQ92").Value = "FP"
$ws.Range("AV92").Value = "TN"
$ws.Range("AR97").Value = "FP"
$ws.Range("AV98").Value = "TN"
$ws.Range("AZ98").Value = "FP"
$ws.Range("AV100").Value = "TN"
$ws.Range("AR102").Value = "TN"
$ws.Range("AV102").Value = "TN"
$ws.Range("AZ103").Value = "TN"
$ws.Range("AR105").Value = "TN"
$ws.Range("AQ106").Value = "FP"
$ws.Range("AZ106").Value = "FP"
$ws.Range("AR107").Value = "FN"
$ws.Range("AR108").Value = "FN"
$ws.Range("AZ112").Value = "FP"
$ws.Range("AQ114").Value = "FP"
$ws.Range("AZ115").Value = "FP"
$ws.Range("AQ116").Value = "FP"
$ws.Range("AV120").Value = "FP"
$ws.Range("AZ120").Value = "TN"
$ws.Range("AZ121").Value = "TN"
$ws.Range("AZ124").Value = "FP"
$ws.Range("AQ125").Value = "TN"
$ws.Range("AR126").Value = "TN"
$ws.Range("AR129").Value = "TN"
$ws.Range("AV129").Value = "TN"
$ws.Range("AZ129").Value = "FP"
$ws.Range("AR130").Value = "TN"
$ws.Range("AV130").Value = "TN"
$ws.Range("AR132").Value = "FP"
$ws.Range("AQ133").Value = "FP"
$ws.Range("AZ133").Value = "FP"
$ws.Range("AV135").Value = "TN"
$ws.Range("AZ135").Value = "FP"
$ws.Range("AR136").Value = "TN"
$ws.Range("AZ136").Value = "TN"
$ws.Range("AR137").Value = "TN"
$ws.Range("AV137").Value = "TN"
$ws.Range("AV139").Value = "TN"
$ws.Range("AZ139").Value = "TN"
$ws.Range("AQ141").Value = "FP"
$ws.Range("AR141").Value = "TN"
$ws.Range("AR143").Value = "TN"
$ws.Range("AZ143").Value = "FP"
$ws.Range("AR146").Value = "TN"
$ws.Range("AR149").Value = "FP"
$ws.Range("AQ152").Value = "FP"
$ws.Range("AZ152").Value = "TN"
$ws.Range("AR154").Value = "FP"
$ws.Range("AQ155").Value = "FP"
$ws.Range("AV156").Value = "TN"
$ws.Range("AR157").Value = "TN"
$ws.Range("AV163").Value = "TN"
$ws.Range("AZ164").Value = "FP"
$ws.Range("AQ165").Value = "FP"
$ws.Range("AR168").Value = "TN"
$ws.Range("AV168").Value = "TN"
$ws.Range("AZ169").Value = "FP"
$ws.Range("AZ170").Value = "FP"
$ws.Range("AV172").Value = "TN"
$ws.Range("AZ174").Value = "TP"
$ws.Range("AR177").Value = "TN"
$ws.Range("AZ178").Value = "FP"
$ws.Range("AZ179").Value = "FP"
$ws.Range("AR181").Value = "TN"
$ws.Range("AQ182").Value = "FP"
$ws.Range("AV184").Value = "TN"
$ws.Range("AV186").Value = "TN"
$ws.Range("AZ187").Value = "FP"
$ws.Range("AV188").Value = "TN"
$ws.Range("AZ189").Value = "FP"
$ws.Range("AZ190").Value = "FP"
$ws.Range("AR192").Value = "TN"
$ws.Range("AR193").Value = "TN"
$ws.Range("AQ196").Value = "FP"
$ws.Range("AZ199").Value = "FP"
$ws.Range("AR200").Value = "TN"
$ws.Range("AR201").Value = "FP"
$ws.Range("AV201").Value = "TN"
